$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.735.14'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '3.493.04'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.39'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.36'
$ws.Range('E6').Value = '  +2.86%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.580'
$ws.Range('E8').Value = '  -1.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.131'
$ws.Range('E9').Value = '  +3.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.12'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '4.100.22'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.28'
$ws.Range('E14').Value = '  +4.86%  '
$ws.Range('D15').Value = '66.812.89'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('D17').Value = '3.517.68'
$ws.Range('E17').Value = '  +0.80%  '
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.30'
$ws.Range('E19').Value = '  +2.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '389.60'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.91'
$ws.Range('E21').Value = '  +0.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.27'
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('E25').Value = '  -1.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000120'
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.11'
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.179'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('E30').Value = '  -2.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.41'
$ws.Range('E31').Value = '  -1.85%  '
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.61'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.35'
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '163.84'
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.874'
$ws.Range('E37').Value = '  -2.09%  '
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.82'
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.61'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('D41').Value = '2.829.98'
$ws.Range('E41').Value = '  +1.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '27.10'
$ws.Range('E42').Value = '  +2.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0731'
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '25.90'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.48'
$ws.Range('E45').Value = '  -0.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.53'
$ws.Range('E46').Value = '  +0.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0299'
$ws.Range('E47').Value = '  -2.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '339.47'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.18'
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.840'
$ws.Range('E51').Value = '  -1.35%  '

Write-Host "Updated cryptos list"